$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "By Crivelatti" -> "By Crivel" + "4" + "tti" (three separate runs,
# each keeping the original italic Arial formatting of the "By Crivelatti"
# run).
# ---------------------------------------------------------------------------
$range = $d.Content
$found = $range.Find.Execute("By Crivelatti")
if ($found) {
    # Shrink the found range down to just "By Crivel" - this keeps the
    # original run's formatting/rPr intact for the first chunk.
    $range.Text = "By Crivel"

    # Insert the "4" chunk right after, then explicitly stamp it with the
    # same Arial italic formatting (new runs from InsertAfter start with no
    # rPr at all).
    $r2 = $d.Range($range.End, $range.End)
    $r2.InsertAfter("4")
    $r2.Font.Name = "Arial"
    $r2.Font.NameBi = "Arial"
    $r2.Font.Italic = $true
    $r2.Font.ItalicBi = $true

    # Insert the "tti" chunk right after that, same formatting treatment.
    $r3 = $d.Range($r2.End, $r2.End)
    $r3.InsertAfter("tti")
    $r3.Font.Name = "Arial"
    $r3.Font.NameBi = "Arial"
    $r3.Font.Italic = $true
    $r3.Font.ItalicBi = $true
}

# ---------------------------------------------------------------------------
# Edit 2: Merge "She walks " / "in, lighting" / " up the dark." (split across
# three runs around a pair of gramStart/gramEnd proofErr markers) into a
# single run "She walks in, lighting up the dark." with no proofErr markers.
# ---------------------------------------------------------------------------
$range2 = $d.Content
$found2 = $range2.Find.Execute("She walks in, lighting up the dark.", $false, $false, $false, $false, $false, $true, 1, $false, "She walks in, lighting up the dark.", 2)
Write-Host "Edit 2 applied:" $found2
